$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.515.32'
$ws.Range("E2").Value = '  +1.77%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.287.48'
$ws.Range("E3").Value = '  +1.05%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.53'
$ws.Range("E5").Value = '  +1.18%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.45'
$ws.Range("E6").Value = '  +7.00%  '

$ws.Range("E7").Value = '  +0.56%  '

$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.497'
$ws.Range("E9").Value = '  +3.57%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.32'
$ws.Range("E10").Value = '  +12.36%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0806'
$ws.Range("E11").Value = '  +1.18%  '

$ws.Range("E12").Value = '  -1.77%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.75'
$ws.Range("E13").Value = '  +2.13%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.646.25'
$ws.Range("E14").Value = '  +1.23%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.59'
$ws.Range("E15").Value = '  +2.77%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.284.84'
$ws.Range("E16").Value = '  +0.57%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.806'
$ws.Range("E17").Value = '  +5.53%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.441.91'
$ws.Range("E18").Value = '  +1.86%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.73'
$ws.Range("E19").Value = '  +1.67%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0917'
$ws.Range("E20").Value = '  +1.58%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.04'
$ws.Range("E21").Value = '  +2.08%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.99'
$ws.Range("E22").Value = '  +1.69%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '243.19'
$ws.Range("E23").Value = '  +1.51%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.61'
$ws.Range("E24").Value = '  +0.97%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.96'
$ws.Range("E25").Value = '  +2.02%  '

$ws.Range("E26").Value = '  +0.05%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.99'
$ws.Range("E27").Value = '  +0.00%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '38.20'
$ws.Range("E28").Value = '  +11.09%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.61'
$ws.Range("E29").Value = '  +0.91%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.12'
$ws.Range("E30").Value = '  +2.16%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '161.11'
$ws.Range("E31").Value = '  +0.34%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.30'
$ws.Range("E32").Value = '  +0.55%  '

$ws.Range("E33").Value = '  +0.11%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.14'
$ws.Range("E34").Value = '  +4.69%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0751'
$ws.Range("E35").Value = '  +0.98%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.38'
$ws.Range("E36").Value = '  +2.93%  '

$ws.Range("E37").Value = '  +2.91%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.87'
$ws.Range("E38").Value = '  +4.30%  '

$ws.Range("E39").Value = '  +0.34%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.115'
$ws.Range("E40").Value = '  -0.56%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.19'
$ws.Range("E41").Value = '  +6.14%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.41'
$ws.Range("E42").Value = '  +15.50%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.005.69'
$ws.Range("E43").Value = '  -0.98%  '

$ws.Range("E44").Value = '  +0.83%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0288'
$ws.Range("E45").Value = '  +3.51%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.05'
$ws.Range("E46").Value = '  +5.71%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.15'
$ws.Range("E47").Value = '  -1.92%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.78'
$ws.Range("E48").Value = '  +3.71%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.54'
$ws.Range("E49").Value = '  +1.36%  '

# Row 50/51 swap (Aave and BitcoinSV changed order)
$ws.Range("D50").NumberFormat = "@"
$ws.Range("B50").Value = 'BitcoinSV'
$ws.Range("C50").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D50").Value = '72.77'
$ws.Range("E50").Value = '  +0.23%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '93.33'
$ws.Range("E51").Value = '  +2.51%  '
